$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the accented name: "André" -> "Andre"
$ws.Range("B2").Value = "Andre"

# Move the active selection from B11 to B2
$ws.Range("B2").Select()
